# moh_515_post_outbreak.xlsx - "worked on household indicators"
#
# 1) Re-case several question labels (sentence-case -> title-case) and
#    rename the first group's "name" column value.
# 2) Add a brand-new "household_indicators" group (7 survey rows) right
#    after the existing group's "end group" row.
# 3) Widen column C (and give column B an explicit width) on the survey
#    sheet to fit the new, longer labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- 1) Title-case the existing labels -------------------------------
$ws.Cells.Item(2, 2).Value = "form_summary"

$ws.Cells.Item(3, 3).Value = "What Is Your Name?"
$ws.Cells.Item(4, 3).Value = "What Is Your Area?"
$ws.Cells.Item(5, 3).Value = "What Is Your Linked Facility?"
$ws.Cells.Item(6, 3).Value = "What County Do You Belong To?"
$ws.Cells.Item(7, 3).Value = "How Many CHPs Are In Your Area?"
$ws.Cells.Item(8, 3).Value = "How Many CHPs Submitted Monthly Report?"

# --- 2) New "household_indicators" group, rows 10-16 ------------------
# row 10: begin group
$ws.Cells.Item(10, 1).Value = "begin group"
$ws.Cells.Item(10, 2).Value = "household_indicators"
$ws.Cells.Item(10, 3).Value = "Household Indicators"

# row 11: total_households
$ws.Cells.Item(11, 1).Value = "integer"
$ws.Cells.Item(11, 2).Value = "total_households"
$ws.Cells.Item(11, 3).Value = "Total Households In The Area?"
$ws.Cells.Item(11, 4).Value = "yes"
$ws.Cells.Item(11, 6).Value = "numbers"

# row 12: new_households
$ws.Cells.Item(12, 1).Value = "integer"
$ws.Cells.Item(12, 2).Value = "new_households"
$ws.Cells.Item(12, 3).Value = "Number Of New Households Registered This Month?"
$ws.Cells.Item(12, 4).Value = "yes"
$ws.Cells.Item(12, 6).Value = "numbers"

# row 13: new_households_visited
$ws.Cells.Item(13, 1).Value = "integer"
$ws.Cells.Item(13, 2).Value = "new_households_visited"
$ws.Cells.Item(13, 3).Value = "Number Of New Households Visited This Month?"
$ws.Cells.Item(13, 4).Value = "yes"
$ws.Cells.Item(13, 6).Value = "numbers"

# row 14: new_households_with_clean_water
$ws.Cells.Item(14, 1).Value = "integer"
$ws.Cells.Item(14, 2).Value = "new_households_with_clean_water"
$ws.Cells.Item(14, 3).Value = "Number Of New Households Visited This Month With Clean Water Access?"
$ws.Cells.Item(14, 4).Value = "yes"

# row 15: new_households_with_latrines
$ws.Cells.Item(15, 1).Value = "integer"
$ws.Cells.Item(15, 2).Value = "new_households_with_latrines"
$ws.Cells.Item(15, 3).Value = "Number Of New Households Visited This Month With Latrines/Toilets?"
$ws.Cells.Item(15, 4).Value = "yes"

# row 16: end group
$ws.Cells.Item(16, 1).Value = "end group"

# Copy the existing body-row style (s="1") onto exactly the populated
# cells of the new rows, without disturbing the shared style table.
$ws.Range("A2").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)
$ws.Range("A11:D11").PasteSpecial(-4122)
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("A12:D12").PasteSpecial(-4122)
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("A13:D13").PasteSpecial(-4122)
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("A14:D14").PasteSpecial(-4122)
$ws.Range("A15:D15").PasteSpecial(-4122)
$ws.Range("A16").PasteSpecial(-4122)

# C15 carries the same highlighted style used elsewhere in the sheet
# (e.g. H8) rather than the plain one.
$ws.Range("H8").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3) Column widths on the survey sheet ------------------------------
$ws.Columns.Item(2).ColumnWidth = 18
$ws.Columns.Item(3).ColumnWidth = 59.33
